# Fix the font size and alignment in the feature comparison table.
#
# This script rewrites the table's paragraph/run formatting by operating on
# the document's OOXML (via Range.WordOpenXML / Range.InsertXML):
#   1. Remove the now-unused <w:bottom w:w="0" w:type="dxa"/> cell margin.
#   2. In every table-cell paragraph, drop the explicit <w:spacing w:after="0"/>
#      and add a paragraph-mark run (<w:rPr>) sized to 21/21 (szCs).
#   3. Center the header cell's paragraph.
#   4. Bump every run from sz=19 to sz=21 and make sure sz=21 runs also carry
#      szCs=21.

$d = $word.ActiveDocument
$range = $d.Content
$xml = $range.WordOpenXML

# 1) Table cell margins: drop the explicit bottom margin of 0.
$xml = $xml -replace '<w:bottom w:w="0" w:type="dxa"/>', ''

# 2) Header cell (unique paragraph / run -- contains the table title text).
#    Handle it on its own because it both centers the paragraph and already
#    had sz=21 (so it only needs szCs added), unlike every other cell.
$oldHeaderPara = '<w:pPr><w:spacing w:after="0"/><w:ind w:left="4"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/></w:rPr><w:t>INCLUDED IN SHAREPOINT ONLINE TO GOOGLE MYDRIVE MIGRATION FEATURES</w:t></w:r>'
$newHeaderPara = '<w:pPr><w:ind w:left="4"/><w:jc w:val="center"/><w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>INCLUDED IN SHAREPOINT ONLINE TO GOOGLE MYDRIVE MIGRATION FEATURES</w:t></w:r>'
$xml = $xml.Replace($oldHeaderPara, $newHeaderPara)

# 3) Remaining "label" cells (left column): pPr has <w:ind w:left="4"/> after
#    the spacing element.  Strip the spacing, add a sized paragraph-mark rPr.
$xml = $xml -replace '<w:pPr><w:spacing w:after="0"/><w:ind w:left="4"/></w:pPr>', '<w:pPr><w:ind w:left="4"/><w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr>'

# 4) Remaining "description" cells (right column): pPr has only spacing.
$xml = $xml -replace '<w:pPr><w:spacing w:after="0"/></w:pPr>', '<w:pPr><w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr>'

# 5) Every run sized 19 becomes 21, with matching szCs added.
$xml = $xml -replace '<w:sz w:val="19"/>', '<w:sz w:val="21"/><w:szCs w:val="21"/>'

$range.InsertXML($xml)
